# Daily update at 8 AM UTC
# Appends the next day's row of data (row 22) to the "Wins Over Time" sheet
# and moves the "last row" date-only formatting from the old last row (21)
# to the new last row (22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 was previously the last row and used the date-only format
# (style index 3). Now that it's no longer the last row, it should use
# the regular interior-row datetime format (style index 2).
$ws.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new last row (22) with the next day's values.
$ws.Range("A22").Value = 45971
$ws.Range("B22").Value = 47
$ws.Range("C22").Value = 55
$ws.Range("D22").Value = 54

# The new last row gets the date-only format that row 21 used to have.
$ws.Range("A22").NumberFormat = "YYYY-MM-DD"
